$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the GoalPoseX value (row 12, column B) from 255 to 208
$ws.Range("B12").Value = 208

# Update the active cell selection from B14 to B13
$ws.Range("B13").Select()
